# Correccion a Diebold Mariano y revision de Cap1
# Rewrites the summary table (rows 2-10) with the corrected
# Comparaciones_Significativas / Proporcion_Sig values and fixes the
# row ordering / labels that result from the correction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Modelo), Column B (Comparaciones_Significativas), Column C
# (Proporcion_Sig), Column D (Mejor_N_Calib), Column E (ECRPS_Mejor)

$ws.Range("A2").Value = "MCPS"
$ws.Range("B2").Value = "3/10"
$ws.Range("C2").Value = 76.8
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 0.606093277050508

$ws.Range("A3").Value = "AV-MCPS"
$ws.Range("B3").Value = "2/10"
$ws.Range("C3").Value = 51.2
$ws.Range("D3").Value = 200
$ws.Range("E3").Value = 0.6292401918857261

$ws.Range("A4").Value = "Sieve Bootstrap"
$ws.Range("B4").Value = "1/10"
$ws.Range("C4").Value = 25.6
$ws.Range("D4").Value = 200
$ws.Range("E4").Value = 0.6013424897804682

$ws.Range("A5").Value = "AREPD"
$ws.Range("B5").Value = "0/10"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 60
$ws.Range("E5").Value = 0.6666371841897598

$ws.Range("A6").Value = "Block Bootstrapping"
$ws.Range("B6").Value = "0/10"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 60
$ws.Range("E6").Value = 0.6094402779078255

$ws.Range("A7").Value = "EnCQR-LSTM"
$ws.Range("B7").Value = "0/10"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 100
$ws.Range("E7").Value = 0.8343460419090608

$ws.Range("A8").Value = "DeepAR"
$ws.Range("B8").Value = "0/10"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 200
$ws.Range("E8").Value = 0.5612074862105157

$ws.Range("A9").Value = "LSPMW"
$ws.Range("B9").Value = "0/10"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 20
$ws.Range("E9").Value = 0.6599565455007916

$ws.Range("A10").Value = "LSPM"
$ws.Range("B10").Value = "0/10"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 60
$ws.Range("E10").Value = 0.6441629191416447
